$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header renames
$ws.Range("D1").Value = "Term Date"
$ws.Range("E1").Value = "Term Reason"

# Row 2
$ws.Range("B2").Value = "Provider"
$ws.Range("F2").Value = "Cole Garrett"
$ws.Range("J2").Value = "RCHN & RCSSD"
$ws.Range("K2").Value = "82-1111113"

# Row 3 - swap C3/D3 values, then other field updates
$ws.Range("C3").Value = "Information not found"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "09/01/2025"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = "Cyrus Hendricks, M.D."
$ws.Range("H3").Value = "Internal Medicine 207R00000X"
$ws.Range("J3").Value = "Mercian Medical Group – P04"
$ws.Range("K3").Value = "45-8888885"
$ws.Range("P3").Value = "P04, 1104, 569"
$ws.Range("Q3").Value = "Medicare, Commercial HMO"

# Row 4
$ws.Range("A4").Value = "Add"
$ws.Range("B4").Value = "Primary Practice Location"
$ws.Range("F4").Value = "Paul Mcmallan, MD"
$ws.Range("J4").Value = "HILABS"
$ws.Range("P4").Value = "P01, P03"
$ws.Range("Q4").Value = "Medicare, Medical"
